$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.305.39'
$ws.Range("E2").Value = '  +0.34%  '
$ws.Range("D3").Value = '1.591.20'
$ws.Range("E3").Value = '  +0.58%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.24%  '
$ws.Range("E6").Value = '  +0.20%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("E8").Value = '  +0.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0609'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.35'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0847'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.27%  '
$ws.Range("D12").Value = '1.815.33'
$ws.Range("E12").Value = '  +0.60%  '
$ws.Range("D13").Value = '1.635.31'
$ws.Range("E13").Value = '  +3.49%  '
$ws.Range("E14").Value = '  +0.31%  '
$ws.Range("E15").Value = '  +0.67%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.45'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").Value = '26.322.13'
$ws.Range("E17").Value = '  +0.41%  '
$ws.Range("E18").Value = '  -0.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.48'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '211.95'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.77%  '
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("E22").Value = '  +0.84%  '
$ws.Range("E23").Value = '  +1.57%  '
$ws.Range("E24").Value = '  -2.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.22'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.26%  '
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.04'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.21%  '
$ws.Range("E28").Value = '  -0.46%  '
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("E30").Value = '  -0.50%  '
$ws.Range("E31").Value = '  +0.86%  '
$ws.Range("E32").Value = '  +0.10%  '
$ws.Range("E33").Value = '  +0.99%  '
$ws.Range("D34").Value = '1.335.55'
$ws.Range("E34").Value = '  +4.10%  '
$ws.Range("E35").Value = '  -0.90%  '
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("E37").Value = '  +0.27%  '
$ws.Range("E38").Value = '  +0.21%  '
$ws.Range("E39").Value = '  -15.21%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.818'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.63%  '
$ws.Range("E41").Value = '  +3.77%  '
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("E43").Value = '  +0.59%  '
$ws.Range("E44").Value = '  -0.70%  '
$ws.Range("D45").Value = '1.727.19'
$ws.Range("E45").Value = '  +0.52%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '61.87'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '87.90'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.80%  '
$ws.Range("E48").Value = '  -3.43%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0981'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.66%  '
$ws.Range("E50").Value = '  -0.68%  '
$ws.Range("E51").Value = '  -0.34%  '
